$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("главная страница")

# Add the four new category rows beneath the existing list, matching the
# bold style already used by the other category cells (A5:A14).
$newValues = @("окружение", "поведение", "ценности", "миссия")
$srcStyle = $ws.Range("A14")

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 15 + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $newValues[$i]
    $srcStyle.Copy()
    $cell.PasteSpecial(-4122) # xlPasteFormats
}

$excel.CutCopyMode = $false

# Update the active selection to reflect the next empty row, as seen after
# the new rows were appended.
$ws.Range("A19").Select()
